# Insert a new row at position 313 (pushes existing rows 313:381 down to 314:382)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("313:313").Insert()

# Populate the newly inserted row 313 with its data (mirrors the other
# rows in this block: same market/region/category/quality/unit/classification)
$ws.Range("A313").Value = 5
$ws.Range("B313").Value = "Macroferia Regional de Talca"
$ws.Range("C313").Value = "Maule"
$ws.Range("D313").Value = 44785
$ws.Range("E313").Value = 7
$ws.Range("F313").Value = 100114013
$ws.Range("G313").Value = "Zanahoria"
$ws.Range("H313").Value = "Sin especificar"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 500
$ws.Range("K313").Value = 9000
$ws.Range("L313").Value = 9000
$ws.Range("M313").Value = 9000
$ws.Range("N313").Value = '$/saco 20 kilos'
$ws.Range("O313").Value = "Región de Ñuble"
$ws.Range("P313").Value = 450
$ws.Range("Q313").Value = 20
$ws.Range("R313").Value = "Hortaliza"
